$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2) | Out-Null
}

# Paragraph 1: authentication side
Replace-Text "password restraints" "password restrictions"
Replace-Text "recovery options, plus the option" "recovery options, as well as the option"
Replace-Text "security concern while" "security concerns while"

# Paragraph 2: location API / TripAdvisor / GooglePlaces
# This paragraph's run has no leading tab, so routing the edit through
# Find.Execute causes the serializer to drop the (semantically irrelevant
# here, but diff-preserved) xml:space="preserve" attribute on <w:t>.
# Assigning Range.Text (using the *full* paragraph range, pilcrow included,
# but with the pilcrow stripped from the replacement string) keeps the
# run/paragraph formatting and the attribute intact, and avoids inserting a
# spurious extra paragraph.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*location API needed*") {
        $trimmed = $t.Substring(0, $t.Length - 1)
        $newText = $trimmed
        $newText = $newText -replace "TripAdvisor to GooglePlaces\.", "TripAdvisor to Google Places."
        $newText = $newText -replace "the GooglePlaces API was due to ongoing issues working with TripAdvisor support\.", "the Google Places API was due to ongoing issues with TripAdvisor support."
        $newText = $newText -replace "Because GooglePlaces provides", "Because Google Places provides"
        $newText = $newText -replace "saving the GooglePlaces ID", "saving the Google Places ID"
        $p.Range.Text = $newText
    }
}

# Paragraph 3: redirect configuration were -> was fixed
Replace-Text "redirect configuration were fixed." "redirect configuration was fixed."

# Paragraph 4: Home Page as well as a Create Itinerary and List Itinerary page.
Replace-Text "Home Page as well as a Create Itinerary and List Itinerary page." "Home Page, as well as Create Itinerary and List Itinerary pages."
